$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the blank separator rows that precede each section header
#    (matches the style already used between the "Devops" and
#    "Python" sections). Insert top-down so each target row index is
#    valid at the moment of insertion.
# ------------------------------------------------------------------
$ws.Rows("12").Insert() | Out-Null
$ws.Rows("22").Insert() | Out-Null
$ws.Rows("27").Insert() | Out-Null
$ws.Rows("32").Insert() | Out-Null

# ------------------------------------------------------------------
# 2. Fill in the "actual time" (column D) numbers that were recorded
#    for the Devops section.
# ------------------------------------------------------------------
$ws.Range("D2").Formula = "=SUM(D3)"
$ws.Range("D3").Value = 0.5

$ws.Range("D5").Formula = "=SUM(D6:D11)"
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 1.5
$ws.Range("D8").Value = 4
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1.5
$ws.Range("D11").Value = 3.5

# ------------------------------------------------------------------
# 3. Shade the section-total cells (columns C & D) with the light
#    gray fill used to highlight subtotal rows.
# ------------------------------------------------------------------
$totalRows = 2, 5, 13, 23, 28, 33
foreach ($r in $totalRows) {
    $ws.Range("C$r`:D$r").Interior.Color = 12566463
}

# ------------------------------------------------------------------
# 4. Row heights / view bits touched by the re-save.
# ------------------------------------------------------------------
$ws.Rows("13").RowHeight = 28.8
$ws.Range("E8").Select() | Out-Null
